$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210, shifting existing rows 210..268 down to 211..269.
$ws.Rows.Item(210).Insert()

# Populate the newly-inserted row 210 with the new price record.
$ws.Range("A210").Value = 10
$ws.Range("B210").Value = "Vega Modelo de Temuco"
$ws.Range("C210").Value = "La Araucanía"
$ws.Range("D210").Value = 44551
$ws.Range("E210").Value = 9
$ws.Range("F210").Value = "Fruta"
$ws.Range("G210").Value = 100108
$ws.Range("H210").Value = "Tropicales y subtropicales"
$ws.Range("I210").Value = 100108002
$ws.Range("J210").Value = "Mango"
$ws.Range("K210").Value = "Sin especificar"
$ws.Range("L210").Value = "Primera"
$ws.Range("M210").Value = 380
$ws.Range("N210").Value = 7000
$ws.Range("O210").Value = 7000
$ws.Range("P210").Value = 7000
$ws.Range("Q210").Value = "$/bandeja 4 kilos"
$ws.Range("R210").Value = "Brasil"
$ws.Range("S210").Value = 1750
$ws.Range("T210").Value = 4
